$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 22.98000000000015
$ws.Range("G2").Value = [double]"1.617237566087226e-09"
$ws.Range("H2").Value = [double]"2.654740167533263e-08"
$ws.Range("K2").Value = 5.949550844899699
$ws.Range("L2").Value = "[3.6492973035673586, 8.24980438623204]"
$ws.Range("M2").Value = [double]"5.537856111459405e-07"
$ws.Range("N2").Value = [double]"5.537856111459405e-07"
$ws.Range("O2").Value = -1.207579158136618
$ws.Range("P2").Value = "[-1.610105544182157, -0.8050527720910781]"
$ws.Range("Q2").Value = [double]"7.528792922428806e-09"
$ws.Range("R2").Value = [double]"1.505758584485761e-08"
$ws.Range("S2").Value = 10.61161639255713
$ws.Range("T2").Value = "[9.337348694200646, 11.885884090913606]"
$ws.Range("W2").Value = 4.416576576576606
$ws.Range("X2").Value = 2.944384384384402
$ws.Range("Y2").Value = 5.888768768768809

# Row 3 updates
$ws.Range("G3").Value = [double]"1.736207622116126e-09"
$ws.Range("H3").Value = [double]"2.654740167533263e-08"
$ws.Range("K3").Value = 5.805253417197701
$ws.Range("L3").Value = "[3.6973101559739234, 7.913196678421478]"
$ws.Range("M3").Value = [double]"9.915286058870265e-08"
$ws.Range("N3").Value = [double]"1.983057211774053e-07"
$ws.Range("O3").Value = -0.1886842434588463
$ws.Range("P3").Value = "[-0.6163685286322318, 0.23900004171453926]"
$ws.Range("Q3").Value = 0.3864383678631347
$ws.Range("R3").Value = 0.3864383678631347
$ws.Range("S3").Value = 10.17317660988999
$ws.Range("T3").Value = "[8.932285947146102, 11.414067272633883]"
$ws.Range("W3").Value = 0.7009009009009048
$ws.Range("X3").Value = -0.8878078078078193
$ws.Range("Y3").Value = 2.289609609609629
